$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b87c9a916a9f4b8894c8333c4d778421530d9e9e/e2e/1a5f1d85-17c5-48e6-b75d-e7df95b95412.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcc1bb5f65ea034d1567865a57e1cb53dcb0a50f/e2e/1a5f1d85-17c5-48e6-b75d-e7df95b95412.md."

# --- zh-cn sheet: row 7 gets a (late/out-of-date) handback recorded ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhHandoffFile = $wsZh.Range("G7").Value
$zhTargetMd = $wsZh.Range("A7").Value

$wsZh.Range("J7").Value = $zhHandoffFile
$wsZh.Range("K7").Value = "2016-08-30 02:58:26"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fcc1bb5f65ea034d1567865a57e1cb53dcb0a50f/e2e/1a5f1d85-17c5-48e6-b75d-e7df95b95412.md", [Type]::Missing, [Type]::Missing, $zhTargetMd)

# --- de-de sheet: row 7 gets the same (late/out-of-date) handback recorded ---
$wsDe = $wb.Worksheets.Item("de-de")

$deHandoffFile = $wsDe.Range("G7").Value
$deTargetMd = $wsDe.Range("A7").Value

$wsDe.Range("J7").Value = $deHandoffFile
$wsDe.Range("K7").Value = "2016-08-30 02:58:33"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fcc1bb5f65ea034d1567865a57e1cb53dcb0a50f/e2e/1a5f1d85-17c5-48e6-b75d-e7df95b95412.md", [Type]::Missing, [Type]::Missing, $deTargetMd)

Write-Output "done"
